# Applies the crypto price/volume refresh described in the commit diff.
# Numeric-looking price strings (e.g. "595.43") must stay as literal text
# (matching the original inlineStr cells), so we briefly force a Text
# number format around the assignment and then restore General so the
# cell's display format is unchanged from before the edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.381.49"
$ws.Range("E2").Value = "  +1.06%  "

$ws.Range("D3").Value = "3.748.71"
$ws.Range("E3").Value = "  -0.92%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.43"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.47"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -0.93%  "

$ws.Range("D7").Value = "3.747.67"
$ws.Range("E7").Value = "  -0.88%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("E9").Value = "  -0.80%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -3.07%  "

$ws.Range("E11").Value = "  -0.21%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.447"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -1.19%  "

$ws.Range("E13").Value = "  -6.63%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.04"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -1.13%  "

$ws.Range("D15").Value = "4.377.56"
$ws.Range("E15").Value = "  -0.93%  "

$ws.Range("D16").Value = "3.745.66"
$ws.Range("E16").Value = "  -0.93%  "

$ws.Range("D17").Value = "68.373.58"
$ws.Range("E17").Value = "  +1.03%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.93"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -3.29%  "

$ws.Range("E19").Value = "  -2.48%  "

$ws.Range("E20").Value = "  -0.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.76"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +2.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "465.93"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -0.57%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.697"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -2.94%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.52"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +0.80%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000145"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -2.52%  "

$ws.Range("E26").Value = "  -0.66%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.98"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -1.54%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.01"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -2.95%  "

$ws.Range("D30").Value = "3.894.98"
$ws.Range("E30").Value = "  -0.97%  "

$ws.Range("E31").Value = "  -4.41%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.31"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -4.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.86"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -2.27%  "

$ws.Range("E34").Value = "  -2.05%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.21"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +0.80%  "

$ws.Range("D37").Value = "3.704.27"
$ws.Range("E37").Value = "  -1.15%  "

$ws.Range("E38").Value = "  -2.64%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.37"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -9.20%  "

$ws.Range("E40").Value = "  +1.08%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -0.45%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.81"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +0.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("B45").Value = "Arweave"
$ws.Range("C45").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "44.18"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +12.50%  "

$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.304"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -2.38%  "

$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.56"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -1.24%  "

$ws.Range("E48").Value = "  -0.93%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "45.95"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +0.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "146.54"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +4.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "388.78"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -1.47%  "

